# Update latest output (run 261)
# - Sheet "Schedule": replace row 2 values and append rows 3-4 (new pump
#   schedule intervals produced by the optimisation run).
# - Sheet "Detailed": a handful of historical/forecast boundary cells
#   shift (Type historical<->forecast, Pump_Status ON<->OFF, revised
#   Price values) and the forecast horizon grows by one day, appending
#   rows 50-97.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Schedule
# ---------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

$scheduleRows = @(
    @(46079.29166666666, 46079.66666666666, 9, 34.02, 261.3783900000001, 7.683080246913583),
    @(46079.91666666666, 46080.125, 5, 18.9, 549.341988, 29.06571365079365),
    @(46080.33333333334, 46080.66666666666, 8, 30.24, 358.5154365, 11.85566919642857)
)

$r = 2
foreach ($row in $scheduleRows) {
    $schedule.Cells.Item($r, 1).Value = $row[0]
    $schedule.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $schedule.Cells.Item($r, 2).Value = $row[1]
    $schedule.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $schedule.Cells.Item($r, 3).Value = $row[2]
    $schedule.Cells.Item($r, 4).Value = $row[3]
    $schedule.Cells.Item($r, 5).Value = $row[4]
    $schedule.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------
# Sheet: Detailed
# ---------------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

# Targeted updates to existing rows (only the cells that actually change)
$detailed.Cells.Item(12, 5).Value = "OFF"
$detailed.Cells.Item(13, 5).Value = "OFF"

$detailed.Cells.Item(14, 2).Value = 61.37567
$detailed.Cells.Item(14, 5).Value = "OFF"

$detailed.Cells.Item(15, 2).Value = 72.91204999999999
$detailed.Cells.Item(15, 5).Value = "OFF"

$detailed.Cells.Item(16, 3).Value = "historical"
$detailed.Cells.Item(17, 3).Value = "historical"
$detailed.Cells.Item(18, 3).Value = "historical"

$detailed.Cells.Item(19, 2).Value = 35.88
$detailed.Cells.Item(19, 3).Value = "historical"

$detailed.Cells.Item(20, 2).Value = 7.89424
$detailed.Cells.Item(20, 3).Value = "historical"

$detailed.Cells.Item(21, 2).Value = 0.51
$detailed.Cells.Item(21, 3).Value = "historical"

$detailed.Cells.Item(22, 2).Value = -0.11174
$detailed.Cells.Item(22, 3).Value = "historical"

$detailed.Cells.Item(23, 2).Value = 0.51
$detailed.Cells.Item(23, 3).Value = "historical"

$detailed.Cells.Item(24, 2).Value = 0.51
$detailed.Cells.Item(24, 3).Value = "historical"

$detailed.Cells.Item(25, 2).Value = 0.69338
$detailed.Cells.Item(25, 3).Value = "historical"

$detailed.Cells.Item(26, 2).Value = 0.5101
$detailed.Cells.Item(26, 3).Value = "historical"

$detailed.Cells.Item(27, 2).Value = 0.51
$detailed.Cells.Item(27, 3).Value = "historical"

$detailed.Cells.Item(28, 2).Value = 14.21228
$detailed.Cells.Item(28, 3).Value = "historical"

$detailed.Cells.Item(29, 2).Value = 0.01063
$detailed.Cells.Item(29, 3).Value = "historical"

$detailed.Cells.Item(30, 2).Value = 0.51
$detailed.Cells.Item(30, 3).Value = "historical"

$detailed.Cells.Item(31, 3).Value = "historical"

$detailed.Cells.Item(32, 2).Value = 19.87439
$detailed.Cells.Item(32, 3).Value = "historical"

$detailed.Cells.Item(33, 2).Value = 17.58706
$detailed.Cells.Item(33, 3).Value = "historical"

$detailed.Cells.Item(34, 2).Value = 0.04374
$detailed.Cells.Item(35, 2).Value = 31.13317
$detailed.Cells.Item(36, 2).Value = 37.89
$detailed.Cells.Item(37, 2).Value = 44.74063
$detailed.Cells.Item(38, 2).Value = 44.63385
$detailed.Cells.Item(39, 2).Value = 64.89
$detailed.Cells.Item(40, 2).Value = 67.24552
$detailed.Cells.Item(41, 2).Value = 78
$detailed.Cells.Item(42, 2).Value = 71.07368
$detailed.Cells.Item(43, 2).Value = 78
$detailed.Cells.Item(44, 2).Value = 78
$detailed.Cells.Item(45, 2).Value = 68.25917

$detailed.Cells.Item(46, 2).Value = 65.01000000000001
$detailed.Cells.Item(46, 5).Value = "ON"

$detailed.Cells.Item(47, 2).Value = 62.07375
$detailed.Cells.Item(47, 5).Value = "ON"

$detailed.Cells.Item(48, 2).Value = 57.06006
$detailed.Cells.Item(48, 5).Value = "ON"

$detailed.Cells.Item(49, 2).Value = 57.06
$detailed.Cells.Item(49, 5).Value = "ON"

# New rows 50-97: the forecast horizon now extends one further day
$newDetailedRows = @(
    @(46080, 59.82237, "forecast", 46080, "ON"),
    @(46080.02083333334, 57.06, "forecast", 46080, "ON"),
    @(46080.04166666666, 56.98, "forecast", 46080, "ON"),
    @(46080.0625, 56.98, "forecast", 46080, "ON"),
    @(46080.08333333334, 37.89, "forecast", 46080, "ON"),
    @(46080.10416666666, 53.4915, "forecast", 46080, "ON"),
    @(46080.125, 52.73522, "forecast", 46080, "OFF"),
    @(46080.14583333334, 56.17892, "forecast", 46080, "OFF"),
    @(46080.16666666666, 56.16364, "forecast", 46080, "OFF"),
    @(46080.1875, 57.06, "forecast", 46080, "OFF"),
    @(46080.20833333334, 57.59255, "forecast", 46080, "OFF"),
    @(46080.22916666666, 59.17295, "forecast", 46080, "OFF"),
    @(46080.25, 65, "forecast", 46080, "OFF"),
    @(46080.27083333334, 78.49348999999999, "forecast", 46080, "OFF"),
    @(46080.29166666666, 76.90311, "forecast", 46080, "OFF"),
    @(46080.3125, 61.91148, "forecast", 46080, "OFF"),
    @(46080.33333333334, 52.1197, "forecast", 46080, "ON"),
    @(46080.35416666666, 47.20723, "forecast", 46080, "ON"),
    @(46080.375, 34.54327, "forecast", 46080, "ON"),
    @(46080.39583333334, 32.94767, "forecast", 46080, "ON"),
    @(46080.41666666666, 30.02988, "forecast", 46080, "ON"),
    @(46080.4375, 18.8444, "forecast", 46080, "ON"),
    @(46080.45833333334, 0.73, "forecast", 46080, "ON"),
    @(46080.47916666666, 34.31134, "forecast", 46080, "ON"),
    @(46080.5, 0.73, "forecast", 46080, "ON"),
    @(46080.52083333334, 0.73, "forecast", 46080, "ON"),
    @(46080.54166666666, 33.83115, "forecast", 46080, "ON"),
    @(46080.5625, 11.92919, "forecast", 46080, "ON"),
    @(46080.58333333334, 0.73, "forecast", 46080, "ON"),
    @(46080.60416666666, 32.41431, "forecast", 46080, "ON"),
    @(46080.625, 0.73, "forecast", 46080, "ON"),
    @(46080.64583333334, 35.88, "forecast", 46080, "ON"),
    @(46080.66666666666, 37.89, "forecast", 46080, "OFF"),
    @(46080.6875, 50.17554, "forecast", 46080, "OFF"),
    @(46080.70833333334, 51.59744, "forecast", 46080, "OFF"),
    @(46080.72916666666, 49.492, "forecast", 46080, "OFF"),
    @(46080.75, 58.60793, "forecast", 46080, "OFF"),
    @(46080.77083333334, 68.95334, "forecast", 46080, "OFF"),
    @(46080.79166666666, 73.00763999999999, "forecast", 46080, "OFF"),
    @(46080.8125, 67.03489999999999, "forecast", 46080, "OFF"),
    @(46080.83333333334, 73.9156, "forecast", 46080, "OFF"),
    @(46080.85416666666, 78, "forecast", 46080, "OFF"),
    @(46080.875, 72.3917, "forecast", 46080, "OFF"),
    @(46080.89583333334, 65, "forecast", 46080, "OFF"),
    @(46080.91666666666, 64.89, "forecast", 46080, "OFF"),
    @(46080.9375, 58.42755, "forecast", 46080, "OFF"),
    @(46080.95833333334, 57.06, "forecast", 46080, "OFF"),
    @(46080.97916666666, 57.06, "forecast", 46080, "OFF")
)

$r = 50
foreach ($row in $newDetailedRows) {
    $detailed.Cells.Item($r, 1).Value = $row[0]
    $detailed.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $detailed.Cells.Item($r, 2).Value = $row[1]
    $detailed.Cells.Item($r, 3).Value = $row[2]
    $detailed.Cells.Item($r, 4).Value = $row[3]
    $detailed.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD"
    $detailed.Cells.Item($r, 5).Value = $row[4]
    $r++
}
